$d = $word.ActiveDocument

# Locate "...Customer Privacy Concerns" inside the "under review" working
# paper entry (Yoo, Hyesung, Song Yao, Luping Sun, Xiaomeng Du, ...), right
# before the closing comma/quote of the title.
$r = $d.Content
$r.Find.Text = "Customer Privacy Concerns"
$r.Find.Forward = $true
$r.Find.Wrap = 1
$found = $r.Find.Execute()

if (-not $found) {
    throw "Could not find target text 'Customer Privacy Concerns'"
}

$insertPos = $r.End
$newTitlePart = " – An Application with Click-stream Data"

# Type the new subtitle text right where the cursor would be, immediately
# before the comma that closes the quoted title.
$ins = $d.Range($insertPos, $insertPos)
$ins.InsertAfter($newTitlePart)
$newCursor = $insertPos + $newTitlePart.Length

# Word leaves its _GoBack edit marker at the point where typing stopped, so
# move it there from wherever it used to be.
$old = $d.Bookmarks("_GoBack")
$old.Delete()
$gobackRange = $d.Range($newCursor, $newCursor)
$d.Bookmarks.Add("_GoBack", $gobackRange)

# The remaining tail (the closing comma, closing curly quote and
# " under review") keeps its own formatting but re-seat it, one mark at a
# time, after the new cursor/bookmark position.
$tailStart = $newCursor
$closeQuote = [char]0x201D
$expectedTail = ",$closeQuote under review"
$tailLen = $expectedTail.Length
$tailRange = $d.Range($tailStart, $tailStart + $tailLen)
if ($tailRange.Text -ne $expectedTail) {
    throw "Unexpected trailing text: [$($tailRange.Text)]"
}
$tailRange.Delete()

$p1 = $d.Range($tailStart, $tailStart)
$p1.InsertAfter(",")
$afterComma = $tailStart + 1
$d.Range($tailStart, $afterComma).Font.NameAscii = "Garamond"
$d.Range($tailStart, $afterComma).Font.Name = "Garamond"

$p2 = $d.Range($afterComma, $afterComma)
$p2.InsertAfter($closeQuote)
$afterQuote = $afterComma + 1
$d.Range($afterComma, $afterQuote).Font.NameAscii = "Garamond"
$d.Range($afterComma, $afterQuote).Font.Name = "Garamond"

$p3 = $d.Range($afterQuote, $afterQuote)
$p3.InsertAfter(" under review")
$afterReview = $afterQuote + (" under review").Length
$d.Range($afterQuote, $afterReview).Font.NameAscii = "Garamond"
$d.Range($afterQuote, $afterReview).Font.Name = "Garamond"

Write-Host "Inserted subtitle and re-seated trailing punctuation + _GoBack bookmark."
